$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tardigrade")

# Update Tardigrade.Framework.AzureStorage release value: 6.1.2 -> 7.0.0
# Also change style from italic to normal (matching style of J8-J11)
$ws.Range("J6").Value = "7.0.0"
$ws.Range("J6").Font.Italic = $false

# Update the "Project / Release" column header J1: 11.5.0 -> 12.0.0
$ws.Range("J1").Value = "12.0.0"

# Update Tardigrade.Framework release value for new column: 9.1.1 -> 10.0.0
$ws.Range("J2").Value = "10.0.0"

# Update Tardigrade.Framework.EntityFrameworkCore release value: 8.3.2 -> 9.0.0
$ws.Range("J8").Value = "9.0.0"

# Update the active selection to J2
$ws.Range("J2").Select()
